$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Majorelle Magdy"
$ws.Range("G3").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G4").Value = "Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G5").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Hanan Ragab, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Veronia Rafat"
$ws.Range("G6").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G7").Value = "Dr. Rana Abo-Zaid, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G8").Value = "Dr. Eman Tantawi, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda"
$ws.Range("G9").Value = "Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Asmaa Reda"
$ws.Range("G10").Value = "Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Sara Wael"
$ws.Range("G11").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Veronia Rafat"
$ws.Range("G13").Value = "Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G14").Value = "Dr. Safa Hany, Dr. Shimaa Ashraf"
$ws.Range("G17").Value = "Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Dina Adel, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen"
$ws.Range("G19").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G23").Value = "Dr. Nourham Mostafa, Dr. Hana Amr"
$ws.Range("G24").Value = "Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Aya Emad, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Remon, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Monica"
$ws.Range("G25").Value = "Dr. Marina Atef, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Remon, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry"
$ws.Range("G27").Value = "Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry"
$ws.Range("G28").Value = "Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Eman Samir Gabry"
$ws.Range("G29").Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Eman Samir Gabry"
$ws.Range("G30").Value = "Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G31").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G32").Value = "Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G33").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Hanan Ragab, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Veronia Rafat"
$ws.Range("G34").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G35").Value = "Dr. Rana Abo-Zaid, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G36").Value = "Dr. Eman Tantawi, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda"
$ws.Range("G37").Value = "Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Asmaa Reda"
$ws.Range("G38").Value = "Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Sara Wael"
$ws.Range("G39").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Veronia Rafat"
$ws.Range("G41").Value = "Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G42").Value = "Dr. Safa Hany, Dr. Shimaa Ashraf"
$ws.Range("G45").Value = "Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Dina Adel, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen"
$ws.Range("G47").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G51").Value = "Dr. Nourham Mostafa, Dr. Hana Amr"
$ws.Range("G52").Value = "Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Aya Emad, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Remon, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Monica"
$ws.Range("G53").Value = "Dr. Marina Atef, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Remon, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry"
$ws.Range("G55").Value = "Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry"
$ws.Range("G56").Value = "Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Eman Samir Gabry"
$ws.Range("G57").Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Eman Samir Gabry"
